$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 327.52
$ws.Range("I15").Value = 327.52
$ws.Range("K15").Value = 982.5599999999999
$ws.Range("M15").Value = -813.5599999999999
$ws.Range("H98").Value = 1155.24
$ws.Range("I98").Value = 1184.8572
$ws.Range("J98").Value = 999.75
$ws.Range("K98").Value = 1184.8572
$ws.Range("L98").Value = 999.75
$ws.Range("M98").Value = 313.1428000000001
$ws.Range("N98").Value = -3995.75
$ws.Range("H122").Value = 1155.24
$ws.Range("I122").Value = 1184.8572
$ws.Range("J122").Value = 999.75
$ws.Range("K122").Value = 3554.5716
$ws.Range("L122").Value = 2999.25
$ws.Range("M122").Value = -1104.5716
$ws.Range("N122").Value = -7899.25
$ws.Range("H129").Value = 941.22986
$ws.Range("I129").Value = 649.1667
$ws.Range("J129").Value = 987.96
$ws.Range("K129").Value = 1947.5001
$ws.Range("L129").Value = 2963.88
$ws.Range("M129").Value = 3052.4999
$ws.Range("N129").Value = -12963.88
$ws.Range("H132").Value = 903.6667
$ws.Range("I132").Value = 696.475
$ws.Range("J132").Value = 2561.2
$ws.Range("K132").Value = 2089.425
$ws.Range("L132").Value = 7683.599999999999
$ws.Range("M132").Value = 440.5749999999998
$ws.Range("N132").Value = -12743.6
$ws.Range("H135").Value = 1240.5172
$ws.Range("I135").Value = 1102.7727
$ws.Range("J135").Value = 1673.4286
$ws.Range("K135").Value = 9924.954299999999
$ws.Range("L135").Value = 15060.8574
$ws.Range("M135").Value = -7389.954299999999
$ws.Range("N135").Value = -20130.8574
$ws.Range("H137").Value = 1319.2653
$ws.Range("I137").Value = 1109.1714
$ws.Range("J137").Value = 1844.5
$ws.Range("K137").Value = 3327.5142
$ws.Range("L137").Value = 5533.5
$ws.Range("M137").Value = -777.5141999999996
$ws.Range("N137").Value = -10633.5
$ws.Range("H138").Value = 2490.9143
$ws.Range("I138").Value = 974.7111
$ws.Range("J138").Value = 5220.08
$ws.Range("K138").Value = 2924.1333
$ws.Range("L138").Value = 15660.24
$ws.Range("M138").Value = 2215.8667
$ws.Range("N138").Value = -25940.24
$ws.Range("H141").Value = 1368.566
$ws.Range("I141").Value = 954.7727
$ws.Range("J141").Value = 3391.5557
$ws.Range("K141").Value = 2864.3181
$ws.Range("L141").Value = 10174.6671
$ws.Range("M141").Value = 2315.6819
$ws.Range("N141").Value = -20534.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5719.774
$ws.Range("I32").Value = 3853.4307
$ws.Range("J32").Value = 16917.834
$ws.Range("K32").Value = 3853.4307
$ws.Range("L32").Value = 16917.834
$ws.Range("M32").Value = -3566.4307
$ws.Range("N32").Value = -17491.834
$ws.Range("H61").Value = 2955.9456
$ws.Range("I61").Value = 3242.2
$ws.Range("J61").Value = 1667.8
$ws.Range("K61").Value = 3242.2
$ws.Range("L61").Value = 1667.8
$ws.Range("M61").Value = -3030.2
$ws.Range("N61").Value = -2091.8
$ws.Range("H74").Value = 1319.0883
$ws.Range("I74").Value = 1374.2778
$ws.Range("J74").Value = 1257
$ws.Range("K74").Value = 1374.2778
$ws.Range("L74").Value = 1257
$ws.Range("M74").Value = -500.2778000000001
$ws.Range("N74").Value = -3005
$ws.Range("H77").Value = 1319.0883
$ws.Range("I77").Value = 1374.2778
$ws.Range("J77").Value = 1257
$ws.Range("K77").Value = 6871.389
$ws.Range("L77").Value = 6285
$ws.Range("M77").Value = -2503.389
$ws.Range("N77").Value = -15021
$ws.Range("H132").Value = 2441805
$ws.Range("I132").Value = 2287.96
$ws.Range("J132").Value = 6253550.5
$ws.Range("K132").Value = 6863.88
$ws.Range("L132").Value = 18760651.5
$ws.Range("M132").Value = -4333.88
$ws.Range("N132").Value = -18765711.5
$ws.Range("H133").Value = 49630.5
$ws.Range("J133").Value = 49630.5
$ws.Range("L133").Value = 49630.5
$ws.Range("N133").Value = -54690.5
$ws.Range("H136").Value = 2955.9456
$ws.Range("I136").Value = 3242.2
$ws.Range("J136").Value = 1667.8
$ws.Range("K136").Value = 9726.599999999999
$ws.Range("L136").Value = 5003.4
$ws.Range("M136").Value = -7176.599999999999
$ws.Range("N136").Value = -10103.4
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4480.972
$ws.Range("I134").Value = 5370.04
$ws.Range("J134").Value = 2460.3635
$ws.Range("K134").Value = 16110.12
$ws.Range("L134").Value = 7381.0905
$ws.Range("M134").Value = -13575.12
$ws.Range("N134").Value = -12451.0905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 205190.1
$ws.Range("I31").Value = 1570.6274
$ws.Range("J31").Value = 604597.5
$ws.Range("K31").Value = 1570.6274
$ws.Range("L31").Value = 604597.5
$ws.Range("M31").Value = -1275.6274
$ws.Range("N31").Value = -605187.5
$ws.Range("H34").Value = 205190.1
$ws.Range("I34").Value = 1570.6274
$ws.Range("J34").Value = 604597.5
$ws.Range("K34").Value = 1570.6274
$ws.Range("L34").Value = 604597.5
$ws.Range("M34").Value = -1368.6274
$ws.Range("N34").Value = -605001.5
$ws.Range("H58").Value = 897.8108
$ws.Range("I58").Value = 550.4039
$ws.Range("J58").Value = 1718.9546
$ws.Range("K58").Value = 550.4039
$ws.Range("L58").Value = 1718.9546
$ws.Range("M58").Value = -347.4039
$ws.Range("N58").Value = -2124.9546
$ws.Range("H132").Value = 1558.2
$ws.Range("I132").Value = 1273.8776
$ws.Range("J132").Value = 2428.9375
$ws.Range("K132").Value = 3821.6328
$ws.Range("L132").Value = 7286.8125
$ws.Range("M132").Value = -1291.6328
$ws.Range("N132").Value = -12346.8125
$ws.Range("H134").Value = 1657.3881
$ws.Range("I134").Value = 1801.6818
$ws.Range("J134").Value = 1381.3478
$ws.Range("K134").Value = 5405.0454
$ws.Range("L134").Value = 4144.0434
$ws.Range("M134").Value = -2870.0454
$ws.Range("N134").Value = -9214.0434
$ws.Range("H136").Value = 897.8108
$ws.Range("I136").Value = 550.4039
$ws.Range("J136").Value = 1718.9546
$ws.Range("K136").Value = 1651.2117
$ws.Range("L136").Value = 5156.8638
$ws.Range("M136").Value = 898.7882999999999
$ws.Range("N136").Value = -10256.8638
$ws.Range("H137").Value = 27095
$ws.Range("I137").Value = 34000
$ws.Range("J137").Value = 24793.334
$ws.Range("K137").Value = 34000
$ws.Range("L137").Value = 24793.334
$ws.Range("M137").Value = -28900
$ws.Range("N137").Value = -34993.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 557.63635
$ws.Range("I107").Value = 366.85715
$ws.Range("K107").Value = 1100.57145
$ws.Range("M107").Value = 819.4285500000001
$ws.Range("H132").Value = 1962169.6
$ws.Range("J132").Value = 3269790.2
$ws.Range("L132").Value = 29428111.8
$ws.Range("N132").Value = -29433171.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3602979
$ws.Range("I122").Value = 4631581
$ws.Range("K122").Value = 13894743
$ws.Range("M122").Value = -13892293
$ws.Range("H132").Value = 2121.9487
$ws.Range("I132").Value = 1795.3077
$ws.Range("J132").Value = 2775.2307
$ws.Range("K132").Value = 5385.9231
$ws.Range("L132").Value = 8325.6921
$ws.Range("M132").Value = -2855.9231
$ws.Range("N132").Value = -13385.6921
$ws.Range("H138").Value = 45429
$ws.Range("J138").Value = 45429
$ws.Range("L138").Value = 45429
$ws.Range("N138").Value = -55709
$ws.Range("H141").Value = 32369.572
$ws.Range("J141").Value = 32369.572
$ws.Range("L141").Value = 32369.572
$ws.Range("N141").Value = -42729.572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1744.7273
$ws.Range("I61").Value = 1785.875
$ws.Range("J61").Value = 1635
$ws.Range("K61").Value = 1785.875
$ws.Range("L61").Value = 1635
$ws.Range("M61").Value = -1583.875
$ws.Range("N61").Value = -2039
$ws.Range("H113").Value = 1744.7273
$ws.Range("I113").Value = 1785.875
$ws.Range("J113").Value = 1635
$ws.Range("K113").Value = 1785.875
$ws.Range("L113").Value = 1635
$ws.Range("M113").Value = 384.125
$ws.Range("N113").Value = -5975
$ws.Range("H136").Value = 6856.1665
$ws.Range("I136").Value = 5034.028
$ws.Range("J136").Value = 10500.444
$ws.Range("K136").Value = 15102.084
$ws.Range("L136").Value = 31501.332
$ws.Range("M136").Value = -12552.084
$ws.Range("N136").Value = -36601.33199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 32803
$ws.Range("J46").Value = 32803
$ws.Range("L46").Value = 32803
$ws.Range("N46").Value = -33265
$ws.Range("H96").Value = 3441
$ws.Range("J96").Value = 3962.3333
$ws.Range("L96").Value = 3962.3333
$ws.Range("N96").Value = -6708.3333
$ws.Range("H123").Value = 30429
$ws.Range("J123").Value = 30429
$ws.Range("L123").Value = 30429
$ws.Range("N123").Value = -40229
$ws.Range("H134").Value = 32803
$ws.Range("J134").Value = 32803
$ws.Range("L134").Value = 98409
$ws.Range("N134").Value = -103479
